# Applies the numeric-value updates described by the commit diff
# (Kujata_Profits: currentAveragePrice / LevePrice* / LeveProfit* recompute)
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1212.625
$ws.Range("I32").Value = 760.4
$ws.Range("J32").Value = 1966.3334
$ws.Range("K32").Value = 760.4
$ws.Range("L32").Value = 1966.3334
$ws.Range("M32").Value = -434.4
$ws.Range("N32").Value = -2618.3334
$ws.Range("H132").Value = 5446.346
$ws.Range("I132").Value = 4457.8667
$ws.Range("K132").Value = 13373.6001
$ws.Range("M132").Value = -10843.6001
$ws.Range("H137").Value = 1483.919
$ws.Range("I137").Value = 1049.8125
$ws.Range("J137").Value = 1814.6666
$ws.Range("K137").Value = 3149.4375
$ws.Range("L137").Value = 5443.9998
$ws.Range("M137").Value = -599.4375
$ws.Range("N137").Value = -10543.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5813.8667
$ws.Range("I32").Value = 5945.3794
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 5945.3794
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -5658.3794
$ws.Range("N32").Value = -2574
$ws.Range("H74").Value = 953.82355
$ws.Range("I74").Value = 481
$ws.Range("K74").Value = 481
$ws.Range("M74").Value = 393
$ws.Range("H77").Value = 953.82355
$ws.Range("I77").Value = 481
$ws.Range("K77").Value = 2405
$ws.Range("M77").Value = 1963
$ws.Range("H119").Value = 28698
$ws.Range("J119").Value = 28698
$ws.Range("L119").Value = 28698
$ws.Range("N119").Value = -38374
$ws.Range("H121").Value = 37999
$ws.Range("J121").Value = 37999
$ws.Range("L121").Value = 37999
$ws.Range("N121").Value = -41493
$ws.Range("H132").Value = 2902.3447
$ws.Range("J132").Value = 3928.25
$ws.Range("L132").Value = 11784.75
$ws.Range("N132").Value = -16844.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2920.8
$ws.Range("I134").Value = 667.6486
$ws.Range("K134").Value = 2002.9458
$ws.Range("M134").Value = 532.0542

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1259.6052
$ws.Range("I31").Value = 1082.421
$ws.Range("K31").Value = 1082.421
$ws.Range("M31").Value = -787.421
$ws.Range("H34").Value = 1259.6052
$ws.Range("I34").Value = 1082.421
$ws.Range("K34").Value = 1082.421
$ws.Range("M34").Value = -880.421
$ws.Range("H52").Value = 34763
$ws.Range("I52").Value = 20709
$ws.Range("J52").Value = 41790
$ws.Range("K52").Value = 20709
$ws.Range("L52").Value = 41790
$ws.Range("M52").Value = -20415
$ws.Range("N52").Value = -42378

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4499.5
$ws.Range("I80").Value = 2998
$ws.Range("K80").Value = 8994
$ws.Range("M80").Value = -8058
$ws.Range("H83").Value = 4499.5
$ws.Range("I83").Value = 2998
$ws.Range("K83").Value = 26982
$ws.Range("M83").Value = -22302
$ws.Range("H107").Value = 9519.583000000001
$ws.Range("I107").Value = 382.5
$ws.Range("J107").Value = 14088.125
$ws.Range("K107").Value = 1147.5
$ws.Range("L107").Value = 42264.375
$ws.Range("M107").Value = 772.5
$ws.Range("N107").Value = -46104.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15130.5
$ws.Range("J39").Value = 15130.5
$ws.Range("L39").Value = 15130.5
$ws.Range("N39").Value = -16194.5
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 803.5
$ws.Range("I22").Value = 499.5
$ws.Range("K22").Value = 499.5
$ws.Range("M22").Value = -204.5
$ws.Range("H27").Value = 803.5
$ws.Range("I27").Value = 499.5
$ws.Range("K27").Value = 499.5
$ws.Range("M27").Value = -392.5
$ws.Range("H46").Value = 4927
$ws.Range("J46").Value = 6711.4287
$ws.Range("L46").Value = 6711.4287
$ws.Range("N46").Value = -7087.4287
$ws.Range("H61").Value = 1332.6666
$ws.Range("I61").Value = 1332.6666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1332.6666
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1130.6666
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 1789.1177
$ws.Range("I68").Value = 1788.4375
$ws.Range("K68").Value = 1788.4375
$ws.Range("M68").Value = -1039.4375
$ws.Range("H71").Value = 1789.1177
$ws.Range("I71").Value = 1788.4375
$ws.Range("K71").Value = 8942.1875
$ws.Range("M71").Value = -5198.1875
$ws.Range("H82").Value = 1893.75
$ws.Range("I82").Value = 1735.7142
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1735.7142
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1374.7142
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 1893.75
$ws.Range("I85").Value = 1735.7142
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1735.7142
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -487.7141999999999
$ws.Range("N85").Value = -5496
$ws.Range("H93").Value = 1040
$ws.Range("J93").Value = 1200
$ws.Range("L93").Value = 1200
$ws.Range("N93").Value = -3696
$ws.Range("H100").Value = 1278.6
$ws.Range("I100").Value = 1296.3334
$ws.Range("J100").Value = 1252
$ws.Range("K100").Value = 1296.3334
$ws.Range("L100").Value = 1252
$ws.Range("M100").Value = -755.3334
$ws.Range("N100").Value = -2334
$ws.Range("H113").Value = 1332.6666
$ws.Range("I113").Value = 1332.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1332.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 837.3334
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 14707675
$ws.Range("I122").Value = 20834948
$ws.Range("K122").Value = 62504844
$ws.Range("M122").Value = -62502394
$ws.Range("H123").Value = 40952
$ws.Range("J123").Value = 40952
$ws.Range("L123").Value = 40952
$ws.Range("N123").Value = -50752

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2527
$ws.Range("H46").Value = 29555.8
$ws.Range("J46").Value = 29555.8
$ws.Range("L46").Value = 29555.8
$ws.Range("N46").Value = -30017.8
$ws.Range("H134").Value = 29555.8
$ws.Range("J134").Value = 29555.8
$ws.Range("L134").Value = 88667.39999999999

